$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 100

$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 100

$ws.Range("B4").Value = 12
$ws.Range("C4").Value = 12
$ws.Range("D4").Value = 100
